# Regenerate save_data column G (K) values.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# This updates the "K" column (column G) values for rows 2-38 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 6
    4  = 7
    5  = 2
    6  = 3
    7  = 6
    8  = 1
    9  = 3
    10 = 3
    11 = 5
    12 = 6
    13 = 7
    14 = 3
    15 = 5
    16 = 8
    17 = 7
    18 = 9
    19 = 8
    20 = 7
    21 = 5
    22 = 5
    23 = 8
    24 = 4
    25 = 9
    26 = 6
    27 = 7
    28 = 7
    29 = 5
    30 = 9
    31 = 9
    32 = 6
    33 = 0
    34 = 5
    35 = 7
    36 = 1
    37 = 5
    38 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
